$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.419591
$ws.Range("H2").Value = 4.258773
$ws.Range("I2").Value = 0.001848767113890483
$ws.Range("J2").Value = 0.001848767113890483
$ws.Range("M2").Value = 91.51130433333333
$ws.Range("N2").Value = 274.533913
$ws.Range("O2").Value = 0.9685519820468944
$ws.Range("P2").Value = 0.9685519820468945
$ws.Range("Q2").Value = 129.908624029861
$ws.Range("R2").Value = 1169.177616268749
$ws.Range("S2").Value = 0.001790627052501744
$ws.Range("T2").Value = 0.001790627052501744

# Row 3
$ws.Range("G3").Value = 1.419591
$ws.Range("H3").Value = 4.258773
$ws.Range("I3").Value = 0.001848767113890483
$ws.Range("J3").Value = 0.001848767113890483
$ws.Range("O3").Value = 0.001425786415744213
$ws.Range("P3").Value = 0.001425786415744214
$ws.Range("Q3").Value = 0.191235942792
$ws.Range("R3").Value = 1.721123485128
$ws.Range("S3").Value = 0.000002635947036859686
$ws.Range("T3").Value = 0.000002635947036859687

# Row 4
$ws.Range("G4").Value = 1.419591
$ws.Range("H4").Value = 4.258773
$ws.Range("I4").Value = 0.001848767113890483
$ws.Range("J4").Value = 0.001848767113890483
$ws.Range("M4").Value = 2.836578333333333
$ws.Range("N4").Value = 8.509734999999999
$ws.Range("O4").Value = 0.03002223153736139
$ws.Range("P4").Value = 0.03002223153736139
$ws.Range("Q4").Value = 4.026781072794999
$ws.Range("R4").Value = 36.241029655155
$ws.Range("S4").Value = 0.00005550411435187947
$ws.Range("T4").Value = 0.00005550411435187947

# Row 5
$ws.Range("I5").Value = 0.9578582377148513
$ws.Range("J5").Value = 0.9578582377148513
$ws.Range("M5").Value = 91.51130433333333
$ws.Range("N5").Value = 274.533913
$ws.Range("O5").Value = 0.9685519820468944
$ws.Range("P5").Value = 0.9685519820468945
$ws.Range("Q5").Value = 67306.50104184783
$ws.Range("R5").Value = 605758.5093766305
$ws.Range("S5").Value = 0.9277354946586647
$ws.Range("T5").Value = 0.9277354946586647

# Row 6
$ws.Range("I6").Value = 0.9578582377148513
$ws.Range("J6").Value = 0.9578582377148513
$ws.Range("O6").Value = 0.001425786415744213
$ws.Range("P6").Value = 0.001425786415744214
$ws.Range("S6").Value = 0.001365701263542527
$ws.Range("T6").Value = 0.001365701263542527

# Row 7
$ws.Range("I7").Value = 0.9578582377148513
$ws.Range("J7").Value = 0.9578582377148513
$ws.Range("M7").Value = 2.836578333333333
$ws.Range("N7").Value = 8.509734999999999
$ws.Range("O7").Value = 0.03002223153736139
$ws.Range("P7").Value = 0.03002223153736139
$ws.Range("Q7").Value = 2086.301402199986
$ws.Range("R7").Value = 18776.71261979987
$ws.Range("S7").Value = 0.02875704179264422
$ws.Range("T7").Value = 0.02875704179264422

# Row 8
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.04029299517125823
$ws.Range("J8").Value = 0.04029299517125823
$ws.Range("M8").Value = 91.51130433333333
$ws.Range("N8").Value = 274.533913
$ws.Range("O8").Value = 0.9685519820468944
$ws.Range("P8").Value = 0.9685519820468945
$ws.Range("Q8").Value = 2831.296338739431
$ws.Range("R8").Value = 25481.66704865488
$ws.Range("S8").Value = 0.0390258603357281
$ws.Range("T8").Value = 0.03902586033572811

# Row 9
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.04029299517125823
$ws.Range("J9").Value = 0.04029299517125823
$ws.Range("O9").Value = 0.001425786415744213
$ws.Range("P9").Value = 0.001425786415744214
$ws.Range("S9").Value = 0.00005744920516482717
$ws.Range("T9").Value = 0.00005744920516482718

# Row 10
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.04029299517125823
$ws.Range("J10").Value = 0.04029299517125823
$ws.Range("M10").Value = 2.836578333333333
$ws.Range("N10").Value = 8.509734999999999
$ws.Range("O10").Value = 0.03002223153736139
$ws.Range("P10").Value = 0.03002223153736139
$ws.Range("Q10").Value = 87.76176788454833
$ws.Range("S10").Value = 0.001209685630365299
$ws.Range("T10").Value = 0.001209685630365299

